# "made all phrases into single word for card.json"
# Break the multi-word native-language phrases (and their english/latin-script
# translations) into individual single-word vocabulary entries, adding a new
# row to the russian and italian sheets for the extra word that was packed
# into the old multi-word phrase.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "russian": split "I have a question" / "I have a sister" phrases
# into plain word pairs, and add a 4th row for "brother".
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("russian")

$ws1.Range("A1").Value = "type"
$ws1.Range("B1").Value = "native"
$ws1.Range("C1").Value = "english"
$ws1.Range("D1").Value = "latin-script"

$ws1.Range("A2").Value = "general"
$ws1.Range("B2").Value = "вопрос"
$ws1.Range("C2").Value = "question"
$ws1.Range("D2").Value = "vopros"

$ws1.Range("A3").Value = "family"
$ws1.Range("B3").Value = "сестра"
$ws1.Range("C3").Value = "sister"
$ws1.Range("D3").Value = "sestra"

# New row 4 - keep the same "vertical center" style as the other data rows
# in column B before writing the value (mirrors the B/F column formatting).
$ws1.Range("B4").VerticalAlignment = -4108
$ws1.Range("A4").Value = "family"
$ws1.Range("B4").Value = "брат"
$ws1.Range("C4").Value = "brother"
$ws1.Range("D4").Value = "brat"

# ---------------------------------------------------------------------
# Sheet "chinese": content unchanged, only shared-string reshuffle happens
# naturally because of the other sheets' edits.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("chinese")

$ws2.Range("A1").Value = "type"
$ws2.Range("B1").Value = "native"
$ws2.Range("C1").Value = "english"
$ws2.Range("D1").Value = "latin-script"

$ws2.Range("A2").Value = "family"
$ws2.Range("B2").Value = "父母"
$ws2.Range("C2").Value = "parents"
$ws2.Range("D2").Value = "Fùmǔ"

$ws2.Range("A3").Value = "country"
$ws2.Range("B3").Value = "法国"
$ws2.Range("C3").Value = "France"
$ws2.Range("D3").Value = "Fàguó"

# ---------------------------------------------------------------------
# Sheet "italian": split "Dov'e la toilette" / "grazie per l'informazione"
# phrases into plain word pairs, and add a 4th row for "thank you".
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("italian")

$ws3.Range("A1").Value = "type"
$ws3.Range("B1").Value = "native"
$ws3.Range("C1").Value = "english"
$ws3.Range("D1").Value = "latin-script"

$ws3.Range("A2").Value = "toilet directions"
$ws3.Range("B2").Value = "toilette"
$ws3.Range("C2").Value = "toilet"

$ws3.Range("A3").Value = "pleasantries"
$ws3.Range("B3").Value = "informazione"
$ws3.Range("C3").Value = "information"

# New row 4 - keep the same "vertical center" style as the other data rows
# in column B before writing the value.
$ws3.Range("B4").VerticalAlignment = -4108
$ws3.Range("A4").Value = "pleasantries"
$ws3.Range("B4").Value = "grazie"
$ws3.Range("C4").Value = "thank you"

# ---------------------------------------------------------------------
# Restore / move the selections to match the edited state: chinese -> B3,
# italian -> A13 (no longer the active tab), russian -> C11 (now active).
# Order matters: selecting a range activates its sheet, and the russian
# sheet must be selected last so it ends up as the active tab.
# ---------------------------------------------------------------------
[void]$ws2.Range("B3").Select()
[void]$ws3.Range("A13").Select()
[void]$ws1.Range("C11").Select()
